$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Fr" (Friday) column values for rows 2-6
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 6

# Update the active cell selection
$ws.Range("I5").Select()
